# Append a new data row (row 7) to the "Artfynd" worksheet, mirroring the
# structure of the existing rows (e.g. row 6).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A7").Value = 112207069
$ws.Range("B7").Value = 90021
$ws.Range("C7").Value = "Ovaliderad"
$ws.Range("D7").Value = "LC"
$ws.Range("E7").Value = 6031
$ws.Range("F7").Value = "Blomkålssvamp"
$ws.Range("G7").Value = "Sparassis crispa"
$ws.Range("H7").Value = "(Wulfen:Fr.) Fr."
$ws.Range("I7").Value = "'1"
$ws.Range("J7").Value = "fruktkroppar"
$ws.Range("K7").Value = "'"
$ws.Range("N7").Value = "'"
$ws.Range("P7").Value = "Skogen Vita bergen, Sm"
$ws.Range("Q7").Value = 496938.2644572215
$ws.Range("R7").Value = 6390349.892939959
$ws.Range("S7").Value = 25
$ws.Range("T7").Value = "Jönköping"
$ws.Range("U7").Value = "Eksjö"
$ws.Range("V7").Value = "Småland"
$ws.Range("W7").Value = "Eksjö"
$ws.Range("Y7").Value = "'2023-09-19"
$ws.Range("Z7").Value = "00:00"
$ws.Range("AA7").Value = "'2023-09-19"
$ws.Range("AB7").Value = "00:00"
$ws.Range("AD7").Value = $false
$ws.Range("AE7").Value = $false
$ws.Range("AF7").Value = "'"
$ws.Range("AG7").Value = $false
$ws.Range("AT7").Value = "'"
$ws.Range("AW7").Value = "Anita Lindström Jensen"
$ws.Range("AX7").Value = "Anita Lindström Jensen"
$ws.Range("AY7").Value = "'"
